$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(5).ColumnWidth = 10.875

$ws.Range("A7").Value = 42650.371493055558
$ws.Range("A7").NumberFormat = "m/d/yy h:mm"
$ws.Range("B7").Value = $false
$ws.Range("C7").Value = 9878.08
$ws.Range("D7").Value = 9920.24
$ws.Range("E7").Value = 104.839996
$ws.Range("F7").Value = 103.95
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = -0.85
$ws.Range("I7").Value = $false
